$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.270.59'
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").Value = '1.929.26'
$ws.Range("E3").Value = '  -0.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9946'
$ws.Range("E4").Value = '  -0.54%  '

$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7476'
$ws.Range("E5").Value = '  +3.58%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '249.76'
$ws.Range("E6").Value = '  -1.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9978'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.19'
$ws.Range("E8").Value = '  -2.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3237'
$ws.Range("E9").Value = '  -3.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07116'
$ws.Range("E10").Value = '  -3.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7880'
$ws.Range("E11").Value = '  -3.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08020'
$ws.Range("E12").Value = '  -1.53%  '

$ws.Range("D13").Value = '1.931.26'
$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.391'
$ws.Range("E14").Value = '  -2.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.68'
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.68'
$ws.Range("E16").Value = '  -1.68%  '

$ws.Range("D17").Value = '30.296.92'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '255.95'
$ws.Range("E18").Value = '  +0.73%  '

$ws.Range("E19").Value = '  -3.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.744'
$ws.Range("E20").Value = '  -2.39%  '

$ws.Range("D21").Value = '2.184.24'
$ws.Range("E21").Value = '  -0.52%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9991'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9980'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.818'
$ws.Range("E24").Value = '  -2.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.542'
$ws.Range("E25").Value = '  -4.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.41'
$ws.Range("E26").Value = '  +1.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.12'
$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.294'
$ws.Range("E28").Value = '  -5.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1318'
$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.358'
$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.531'
$ws.Range("E31").Value = '  -2.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.427'
$ws.Range("E32").Value = '  -1.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.152'
$ws.Range("E33").Value = '  -2.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05113'
$ws.Range("E34").Value = '  -4.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.292'
$ws.Range("E35").Value = '  -1.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7500'
$ws.Range("E36").Value = '  -1.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.760'
$ws.Range("E37").Value = '  +0.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01986'
$ws.Range("E38").Value = '  -0.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.795'
$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.38'
$ws.Range("E40").Value = '  -3.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.408'
$ws.Range("E41").Value = '  -3.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4531'
$ws.Range("E42").Value = '  -1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.992'
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8442'
$ws.Range("E44").Value = '  -0.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9985'
$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.59'
$ws.Range("E46").Value = '  -1.62%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.877'
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.509'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '977.57'
$ws.Range("E49").Value = '  +10.99%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.90'
$ws.Range("E50").Value = '  -0.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4213'
$ws.Range("E51").Value = '  -0.17%  '
